$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: result changed from "Não" to "NA", and the observation text was
# rewritten to explain that no audit deadline was specified.
$ws.Range("F9").Value = "Não há prazo para auditoria especificado no Plano de Projeto, impossibilitando determinar se já deveria ter ocorrido auditoria ou não para a baseline criada."
$ws.Range("C9").Value = "NA"

# Selection moved from F9 to C9.
[void]$ws.Range("C9").Select()
